$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 4009.3333
$ws.Range("I38").Value = 102.333336
$ws.Range("J38").Value = 7916.3335
$ws.Range("K38").Value = 307.000008
$ws.Range("L38").Value = 23749.0005
$ws.Range("M38").Value = 64.99999200000002
$ws.Range("N38").Value = -24493.0005
$ws.Range("H39").Value = 382.5
$ws.Range("I39").Value = 550
$ws.Range("J39").Value = 298.75
$ws.Range("K39").Value = 1650
$ws.Range("L39").Value = 896.25
$ws.Range("M39").Value = -1354
$ws.Range("N39").Value = -1488.25
$ws.Range("H51").Value = 8942.5
$ws.Range("J51").Value = 8590
$ws.Range("L51").Value = 8590
$ws.Range("N51").Value = -9558
$ws.Range("H64").Value = 71435176
$ws.Range("I64").Value = 5082.6665
$ws.Range("J64").Value = 125007750
$ws.Range("K64").Value = 5082.6665
$ws.Range("L64").Value = 125007750
$ws.Range("M64").Value = -4834.6665
$ws.Range("N64").Value = -125008246
$ws.Range("H67").Value = 71435176
$ws.Range("I67").Value = 5082.6665
$ws.Range("J67").Value = 125007750
$ws.Range("K67").Value = 5082.6665
$ws.Range("L67").Value = 125007750
$ws.Range("M67").Value = -4224.6665
$ws.Range("N67").Value = -125009466
$ws.Range("H70").Value = 1536.7
$ws.Range("I70").Value = 1678.2
$ws.Range("J70").Value = 1395.2
$ws.Range("K70").Value = 5034.6
$ws.Range("L70").Value = 4185.6
$ws.Range("M70").Value = -4764.6
$ws.Range("N70").Value = -4725.6
$ws.Range("H73").Value = 1536.7
$ws.Range("I73").Value = 1678.2
$ws.Range("J73").Value = 1395.2
$ws.Range("K73").Value = 5034.6
$ws.Range("L73").Value = 4185.6
$ws.Range("M73").Value = -4098.6
$ws.Range("N73").Value = -6057.6
$ws.Range("H93").Value = 50000
$ws.Range("J93").Value = 50000
$ws.Range("L93").Value = 50000
$ws.Range("N93").Value = -54992
$ws.Range("H94").Value = 9998
$ws.Range("I94").Value = 9998
$ws.Range("K94").Value = 9998
$ws.Range("M94").Value = -9547
$ws.Range("H116").Value = 4500
$ws.Range("I116").Value = 3333.3333
$ws.Range("K116").Value = 3333.3333
$ws.Range("M116").Value = 108.6667000000002
$ws.Range("H137").Value = 2506016.8
$ws.Range("I137").Value = 4167730
$ws.Range("K137").Value = 12503190
$ws.Range("M137").Value = -12500640
$ws.Range("H138").Value = 327708
$ws.Range("I138").Value = 3028.7334
$ws.Range("J138").Value = 534950.0600000001
$ws.Range("K138").Value = 9086.200199999999
$ws.Range("L138").Value = 1604850.18
$ws.Range("M138").Value = -3946.200199999999
$ws.Range("N138").Value = -1615130.18

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3232.678
$ws.Range("I32").Value = 3275.9473
$ws.Range("K32").Value = 3275.9473
$ws.Range("M32").Value = -2988.9473
$ws.Range("H61").Value = 4999.9614
$ws.Range("I61").Value = 3443.9443
$ws.Range("J61").Value = 8501
$ws.Range("K61").Value = 3443.9443
$ws.Range("L61").Value = 8501
$ws.Range("M61").Value = -3231.9443
$ws.Range("N61").Value = -8925
$ws.Range("H74").Value = 216966.23
$ws.Range("I74").Value = 558095.3
$ws.Range("J74").Value = 3760.5625
$ws.Range("K74").Value = 558095.3
$ws.Range("L74").Value = 3760.5625
$ws.Range("M74").Value = -557221.3
$ws.Range("N74").Value = -5508.5625
$ws.Range("H77").Value = 216966.23
$ws.Range("I77").Value = 558095.3
$ws.Range("J77").Value = 3760.5625
$ws.Range("K77").Value = 2790476.5
$ws.Range("L77").Value = 18802.8125
$ws.Range("M77").Value = -2786108.5
$ws.Range("N77").Value = -27538.8125
$ws.Range("H136").Value = 4999.9614
$ws.Range("I136").Value = 3443.9443
$ws.Range("J136").Value = 8501
$ws.Range("K136").Value = 10331.8329
$ws.Range("L136").Value = 25503
$ws.Range("M136").Value = -7781.832900000001
$ws.Range("N136").Value = -30603

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2351.4285
$ws.Range("I58").Value = 1832.9615
$ws.Range("K58").Value = 1832.9615
$ws.Range("M58").Value = -1629.9615
$ws.Range("H62").Value = 10009190
$ws.Range("I62").Value = 25004106
$ws.Range("J62").Value = 12580.167
$ws.Range("K62").Value = 25004106
$ws.Range("L62").Value = 12580.167
$ws.Range("M62").Value = -25003482
$ws.Range("N62").Value = -13828.167
$ws.Range("H65").Value = 10009190
$ws.Range("I65").Value = 25004106
$ws.Range("J65").Value = 12580.167
$ws.Range("K65").Value = 125020530
$ws.Range("L65").Value = 62900.835
$ws.Range("M65").Value = -125017410
$ws.Range("N65").Value = -69140.83499999999
$ws.Range("H132").Value = 2084.9565
$ws.Range("I132").Value = 1687.1052
$ws.Range("K132").Value = 5061.3156
$ws.Range("M132").Value = -2531.3156
$ws.Range("H134").Value = 2003.6227
$ws.Range("I134").Value = 1946.4131
$ws.Range("K134").Value = 5839.2393
$ws.Range("M134").Value = -3304.2393
$ws.Range("H136").Value = 2351.4285
$ws.Range("I136").Value = 1832.9615
$ws.Range("K136").Value = 5498.8845
$ws.Range("M136").Value = -2948.8845

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 7269.8667
$ws.Range("J39").Value = 7782
$ws.Range("L39").Value = 23346
$ws.Range("N39").Value = -23934
$ws.Range("H45").Value = 6000
$ws.Range("J45").Value = 6000
$ws.Range("L45").Value = 18000
$ws.Range("N45").Value = -19064
$ws.Range("H109").Value = 2412
$ws.Range("J109").Value = 3966.5
$ws.Range("L109").Value = 11899.5
$ws.Range("N109").Value = -13979.5
$ws.Range("H120").Value = 11234.25
$ws.Range("I120").Value = 10062.333
$ws.Range("J120").Value = 14750
$ws.Range("K120").Value = 30186.999
$ws.Range("L120").Value = 44250
$ws.Range("M120").Value = -25348.999
$ws.Range("N120").Value = -53926
$ws.Range("H122").Value = 1691.8636
$ws.Range("J122").Value = 1898.2354
$ws.Range("L122").Value = 17084.1186
$ws.Range("N122").Value = -21984.1186
$ws.Range("H123").Value = 1335.6428
$ws.Range("I123").Value = 992.2308
$ws.Range("J123").Value = 5800
$ws.Range("K123").Value = 2976.6924
$ws.Range("L123").Value = 17400
$ws.Range("M123").Value = -526.6923999999999
$ws.Range("N123").Value = -22300
$ws.Range("H125").Value = 3666.3333
$ws.Range("J125").Value = 7000
$ws.Range("L125").Value = 21000
$ws.Range("N125").Value = -30840
$ws.Range("H126").Value = 6515
$ws.Range("I126").Value = 6030
$ws.Range("J126").Value = 7000
$ws.Range("K126").Value = 18090
$ws.Range("L126").Value = 21000
$ws.Range("M126").Value = -13150
$ws.Range("N126").Value = -30880
$ws.Range("H131").Value = 13801.333
$ws.Range("I131").Value = 25421.666
$ws.Range("J131").Value = 2181
$ws.Range("K131").Value = 76264.99800000001
$ws.Range("L131").Value = 6543
$ws.Range("M131").Value = -71224.99800000001
$ws.Range("N131").Value = -16623
$ws.Range("H134").Value = 1608.1818
$ws.Range("I134").Value = 1608.1818
$ws.Range("K134").Value = 4824.5454
$ws.Range("M134").Value = 245.4546
$ws.Range("H139").Value = 3078.0625
$ws.Range("I139").Value = 2156.125
$ws.Range("K139").Value = 6468.375
$ws.Range("M139").Value = -1328.375
$ws.Range("H140").Value = 9096.286
$ws.Range("I140").Value = 8642.154
$ws.Range("K140").Value = 25926.462
$ws.Range("M140").Value = -20746.462

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 2493.3635
$ws.Range("I113").Value = 2392.7
$ws.Range("J113").Value = 3500
$ws.Range("K113").Value = 2392.7
$ws.Range("L113").Value = 3500
$ws.Range("M113").Value = -222.6999999999998
$ws.Range("N113").Value = -7840

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H43").Value = 66599.60000000001
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 66599.60000000001
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 66599.60000000001
$ws.Range("M43").ClearContents()
$ws.Range("N43").Value = -66985.60000000001
$ws.Range("H55").Value = 355.7
$ws.Range("I55").Value = 151.83333
$ws.Range("J55").Value = 661.5
$ws.Range("K55").Value = 151.83333
$ws.Range("L55").Value = 661.5
$ws.Range("M55").Value = 21.16667000000001
$ws.Range("N55").Value = -1007.5
$ws.Range("H122").Value = 3065.4666
$ws.Range("I122").Value = 3265.6667
$ws.Range("K122").Value = 9797.000100000001
$ws.Range("M122").Value = -7347.000100000001

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H51").Value = 15500
$ws.Range("J51").Value = 35000
$ws.Range("L51").Value = 35000
$ws.Range("N51").Value = -36020
$ws.Range("H62").Value = 7007.375
$ws.Range("I62").Value = 4289.75
$ws.Range("K62").Value = 4289.75
$ws.Range("M62").Value = -3665.75
$ws.Range("H65").Value = 7007.375
$ws.Range("I65").Value = 4289.75
$ws.Range("K65").Value = 21448.75
$ws.Range("M65").Value = -18328.75
$ws.Range("H81").Value = 3744.2856
$ws.Range("I81").Value = 1666.3636
$ws.Range("K81").Value = 3332.7272
$ws.Range("M81").Value = -2271.7272
$ws.Range("H84").Value = 3744.2856
$ws.Range("I84").Value = 1666.3636
$ws.Range("K84").Value = 16663.636
$ws.Range("M84").Value = -11359.636
$ws.Range("H122").Value = 8930620
$ws.Range("I122").Value = 1968.0952
$ws.Range("K122").Value = 5904.2856
$ws.Range("M122").Value = -3454.2856
$ws.Range("H126").Value = 2032.3334
$ws.Range("I126").Value = 2049
$ws.Range("J126").Value = 1999
$ws.Range("K126").Value = 6147
$ws.Range("L126").Value = 5997
$ws.Range("M126").Value = -3677
$ws.Range("N126").Value = -10937
